$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Addresses already present in rows 2-5 (as shared strings), to be logged
# again (repeated) followed by a brand new address - this mirrors a script
# that appends a console/log entry each time it runs over the address list.
$addresses = @(
    "1 Crossgates Mall Road",
    "Duke Rd & Walden Ave",
    "630 Old Country Rd.",
    "160 Walt Whitman Rd.",
    "1 Crossgates Mall Road",
    "Duke Rd & Walden Ave",
    "630 Old Country Rd.",
    "160 Walt Whitman Rd.",
    "Matrix house Milton Keynes UK"
)

$row = 6
foreach ($addr in $addresses) {
    $ws.Range("A$row").Value = $addr
    $ws.Range("A5").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Rows.Item($row).RowHeight = 13.65
    $row = $row + 1
}

# Trailing blank, but still-styled, log rows (rows 15-24) ready to receive
# future console output lines.
while ($row -le 24) {
    $ws.Range("A5").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Rows.Item($row).RowHeight = 13.65
    $row = $row + 1
}
